# inlineForeignTabs 可编辑聚合表格, isDeleteCascade 级联删除
# Remove the "is_locked_lbl" / "is_enabled_lbl" option columns from the
# options.xlsx template: delete columns D:E (their header cells) and shift
# the remaining columns (order_by, rem) left so the row becomes
# lbl, ky, val, order_by, rem.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1:E1").Delete(-4159)
